$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 297 (shifts existing rows 297.. down by one,
# preserving formatting the same way Excel's native "Insert Row" does).
$ws.Rows.Item(297).Insert()

# Populate the newly inserted row 297 with the new data record.
$ws.Range("A297").Value = 6
$ws.Range("B297").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C297").Value = "Metropolitana"
$ws.Range("D297").Value = 45211
$ws.Range("E297").Value = 13
$ws.Range("F297").Value = 100112001
$ws.Range("G297").Value = "Berenjena"
$ws.Range("H297").Value = "Sin especificar"
$ws.Range("I297").Value = "Primera"
$ws.Range("J297").Value = 450
$ws.Range("K297").Value = 5500
$ws.Range("L297").Value = 6000
$ws.Range("M297").Value = 5778
$ws.Range("N297").Value = "$/caja 50 unidades"
$ws.Range("O297").Value = "Región de Arica y Parinacota"
$ws.Range("P297").Value = 116
$ws.Range("Q297").Value = 50
$ws.Range("R297").Value = "Hortaliza"
